$d = $word.ActiveDocument

# --- Edit 1: add a new bulleted "Things to do" item right after
#     "Break Canvas3D into 2 parts ..." (same numbered list, numId 2) ---
$rng = $d.Content
$found = $rng.Find.Execute("Break Canvas3D into 2 parts nonAWTCanvas3D and normalCanvas3D", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $para = $rng.Paragraphs(1)
    $para.Range.InsertParagraphAfter()
    $newPara = $para.Next()
    $insertPoint = $d.Range($newPara.Range.Start, $newPara.Range.Start)
    $insertPoint.Text = "Now using GL2ES2 profile depth buffer returns 0 (but appears to work)"
}

# --- Edit 2: merge the split "import" + " " runs (with stray gramStart/
#     gramEnd proof-error markers) into a single "import " run, just for
#     the "import java.awt.image.ComponentColorModel;" line ---
$rng2 = $d.Content
$found2 = $rng2.Find.Execute("java.awt.image.ComponentColorModel", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found2) {
    $para2 = $rng2.Paragraphs(1)
    $xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w:rsidR="001C5F65" w:rsidRDefault="001C5F65" w:rsidP="001C5F65"><w:pPr><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:adjustRightInd w:val="0"/><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">import </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>java.awt.image.ComponentColorModel</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>;</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    [void]$para2.Range.InsertXML($xml)
}

Write-Output "done"
